# Apply roster update: remove one row (18 -> 17 data rows) and
# rewrite player/position/team assignments per the new roster.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Original sheet had 18 data rows (rows 2-19); new layout has 17 (rows 2-18).
$ws.Rows.Item(19).Delete()

$ws.Cells.Item(2,1).Value = "Isaiah Collier"
$ws.Cells.Item(2,2).Value = "PG"
$ws.Cells.Item(2,3).Value = "Utah Jazz"
$ws.Cells.Item(3,1).Value = "Tyler Herro"
$ws.Cells.Item(3,2).Value = "PG,SG"
$ws.Cells.Item(3,3).Value = "Miami Heat"
$ws.Cells.Item(4,1).Value = "De'Aaron Fox"
$ws.Cells.Item(4,2).Value = "PG"
$ws.Cells.Item(4,3).Value = "San Antonio Spurs"
$ws.Cells.Item(5,1).Value = "Nick Smith Jr."
$ws.Cells.Item(5,2).Value = "PG,SG"
$ws.Cells.Item(5,3).Value = "Charlotte Hornets"
$ws.Cells.Item(6,1).Value = "Miles Bridges"
$ws.Cells.Item(6,2).Value = "SF,PF"
$ws.Cells.Item(6,3).Value = "Charlotte Hornets"
$ws.Cells.Item(7,1).Value = "Evan Mobley"
$ws.Cells.Item(7,2).Value = "PF,C"
$ws.Cells.Item(7,3).Value = "Cleveland Cavaliers"
$ws.Cells.Item(8,1).Value = "Scottie Barnes"
$ws.Cells.Item(8,2).Value = "PG,SG,SF,PF"
$ws.Cells.Item(8,3).Value = "Toronto Raptors"
$ws.Cells.Item(9,1).Value = "Nicolas Claxton"
$ws.Cells.Item(9,2).Value = "C"
$ws.Cells.Item(9,3).Value = "Brooklyn Nets"
$ws.Cells.Item(10,1).Value = "Brook Lopez"
$ws.Cells.Item(10,2).Value = "C"
$ws.Cells.Item(10,3).Value = "Milwaukee Bucks"
$ws.Cells.Item(11,1).Value = "Nikola Vucevic"
$ws.Cells.Item(11,2).Value = "PF,C"
$ws.Cells.Item(11,3).Value = "Chicago Bulls"
$ws.Cells.Item(12,1).Value = "Shaedon Sharpe"
$ws.Cells.Item(12,2).Value = "SG,SF"
$ws.Cells.Item(12,3).Value = "Portland Trail Blazers"
$ws.Cells.Item(13,1).Value = "DeMar DeRozan"
$ws.Cells.Item(13,2).Value = "SF,PF"
$ws.Cells.Item(13,3).Value = "Sacramento Kings"
$ws.Cells.Item(14,1).Value = "Ja Morant"
$ws.Cells.Item(14,2).Value = "PG"
$ws.Cells.Item(14,3).Value = "Memphis Grizzlies"
$ws.Cells.Item(15,1).Value = "Josh Giddey"
$ws.Cells.Item(15,2).Value = "PG,SG,SF"
$ws.Cells.Item(15,3).Value = "Chicago Bulls"
$ws.Cells.Item(16,1).Value = "Mikal Bridges"
$ws.Cells.Item(16,2).Value = "SG,SF,PF"
$ws.Cells.Item(16,3).Value = "New York Knicks"
$ws.Cells.Item(17,1).Value = "Luka Doncic"
$ws.Cells.Item(17,2).Value = "PG,SG"
$ws.Cells.Item(17,3).Value = "Los Angeles Lakers"
$ws.Cells.Item(18,1).Value = "P.J. Washington"
$ws.Cells.Item(18,2).Value = "SF,PF"
$ws.Cells.Item(18,3).Value = "Dallas Mavericks"
